$d = $word.ActiveDocument

# Pull the canonical package XML (pkg:package wrapper around every part,
# including word/document.xml) so we can splice in the new table rows
# with full control over every element/attribute, exactly as specified
# by the target OOXML diff.
$xml = $d.Content.WordOpenXML

# ---------------------------------------------------------------------
# 1) New "Brand" row - inserted right after the "Campaign Type" row and
#    before the "Business & Marketing Objectives:" row.
# ---------------------------------------------------------------------
$anchorBrand = '{{PLACEHOLDER_CAMPAIGN_TYPE}}</w:t></w:r></w:p></w:tc></w:tr>'
$rowBrand = '<w:tr><w:trPr><w:trHeight w:val="288"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="3166" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>Brand</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="6345" w:type="dxa"/><w:gridSpan w:val="2"/></w:tcPr><w:p><w:r><w:t>{{PLACEHOLDER_BRAND_NAME}}</w:t></w:r></w:p></w:tc></w:tr>'
$xml = $xml.Replace($anchorBrand, $anchorBrand + $rowBrand)

# ---------------------------------------------------------------------
# 2) New "Email Subject Line:" + "Email Content:" rows - inserted right
#    after the "Core Message & Positioning:" row and before the
#    "Creative Assets Required:" row.
# ---------------------------------------------------------------------
$anchorEmail = '{{PLACEHOLDER_CORE_MESSAGE}}</w:t></w:r></w:p></w:tc></w:tr>'
$rowEmailSubject = '<w:tr><w:trPr><w:trHeight w:val="594"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="3166" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>Email Subject Line:</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="6345" w:type="dxa"/><w:gridSpan w:val="2"/><w:tcBorders><w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders></w:tcPr><w:p><w:pPr><w:ind w:right="-104"/></w:pPr><w:r><w:t>{{PLACEHOLDER_EMAIL_SUBJECTLINE}}</w:t></w:r></w:p></w:tc></w:tr>'
$rowEmailContent = '<w:tr><w:trPr><w:trHeight w:val="594"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="3166" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>Email Content:</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="6345" w:type="dxa"/><w:gridSpan w:val="2"/><w:tcBorders><w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders></w:tcPr><w:p><w:pPr><w:ind w:right="-104"/></w:pPr><w:r><w:t>{{PLACEHOLDER_EMAIL_CONTENT}}</w:t></w:r></w:p></w:tc></w:tr>'
$xml = $xml.Replace($anchorEmail, $anchorEmail + $rowEmailSubject + $rowEmailContent)

# ---------------------------------------------------------------------
# 3) New "Comments and Approval:" row - inserted at the very end of the
#    table, right after the "Roles & Responsibilities:" row.
# ---------------------------------------------------------------------
$anchorComments = '{{PLACEHOLDER_ROLES}}</w:t></w:r></w:p></w:tc></w:tr>'
$rowComments = '<w:tr><w:trPr><w:trHeight w:val="288"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="3166" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>Comments and Approval:</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="6345" w:type="dxa"/><w:gridSpan w:val="2"/><w:tcBorders><w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders></w:tcPr><w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="1052"/></w:tabs><w:ind w:right="-104"/></w:pPr></w:p></w:tc></w:tr>'
$xml = $xml.Replace($anchorComments, $anchorComments + $rowComments)

# Push the modified package XML back into the document in one shot.
$d.Content.InsertXML($xml)

Write-Host "Inserted Brand, Email Subject Line, Email Content and Comments and Approval rows"
